$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$r = $ws1.Range("C1")
try {
  $r.BorderAround(1, 2, 1)
  Write-Host "ok"
} catch {
  Write-Host "err: $_"
}
